$wb = $excel.ActiveWorkbook

# sheet1 -> Worksheets.Item(1)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 12892
$ws.Range("G2").Value = 85
$ws.Range("F3").Value = 7218
$ws.Range("G3").Value = 80
$ws.Range("G6").Value = 75
$ws.Range("F10").Value = 1015
$ws.Range("F11").Value = 153
$ws.Range("F13").Value = 1027
$ws.Range("F14").Value = 9
$ws.Range("F16").Value = 1025
$ws.Range("F18").Value = 257
$ws.Range("F20").Value = 25
$ws.Range("F21").Value = 286
$ws.Range("F24").Value = 199
$ws.Range("F25").Value = 382
$ws.Range("F26").Value = 5264
$ws.Range("F28").Value = 1448
$ws.Range("F30").Value = 1721
$ws.Range("F31").Value = 79
$ws.Range("F32").Value = 72
$ws.Range("F33").Value = 1376
$ws.Range("F38").Value = 3744

# sheet3 -> Worksheets.Item(3)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 9294
$ws.Range("F3").Value = 563
$ws.Range("F4").Value = 2042

# sheet4 -> Worksheets.Item(4)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 9294
$ws.Range("F3").Value = 563
$ws.Range("F4").Value = 2042
$ws.Range("F5").Value = 12892
$ws.Range("G5").Value = 85
$ws.Range("F6").Value = 7218
$ws.Range("G6").Value = 80
$ws.Range("G9").Value = 75
$ws.Range("F10").Value = 1015
$ws.Range("F11").Value = 153
$ws.Range("F13").Value = 1027
$ws.Range("F14").Value = 9
$ws.Range("F16").Value = 1025
$ws.Range("F17").Value = 257
$ws.Range("F19").Value = 25
$ws.Range("F20").Value = 286
$ws.Range("F26").Value = 199
$ws.Range("F27").Value = 382
$ws.Range("F28").Value = 5264
$ws.Range("F30").Value = 1448
$ws.Range("F35").Value = 1721
$ws.Range("F36").Value = 79
$ws.Range("F37").Value = 72
$ws.Range("F38").Value = 1376
$ws.Range("F47").Value = 3744
